$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header labels: "<name>_old" -> "<name>_FV2210" and "<name>_new" -> "<name>_FV2304"
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_old$", "_FV2210")
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace "_new$", "_FV2304")
}

# 2) Freeze the header row (pane split after row 1) and keep focus on the frozen pane
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# 3) Turn the data range into an Excel Table ("Table1"), preserving the existing header
#    formatting (bold/shaded/bordered header cells must keep using the original cell
#    style instead of Excel auto-generating a dedicated header dxf for the table).
$header = $ws.Range("A1:U1")
$scratch = $ws.Range("A200:U200")
$header.Copy() | Out-Null
$scratch.PasteSpecial(-4122) | Out-Null          # xlPasteFormats
$header.ClearFormats() | Out-Null

$dataRange = $ws.Range("A1:U90")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

$scratch.Copy() | Out-Null
$header.PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$excel.CutCopyMode = $false

# drop the scratch row entirely so it leaves no structural trace behind
$ws.Rows.Item(200).Delete() | Out-Null

$ws.Range("A1").Select() | Out-Null
